$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Q11) Merge two sorted arrays into one  ->  highlight the whole line yellow
# ---------------------------------------------------------------------------
$q11 = $d.Paragraphs(11)
$q11.Range.HighlightColorIndex = 7   # wdYellow

# ---------------------------------------------------------------------------
# Q13) split "in  a" out into its own run wrapped in gramStart/gramEnd
#      proofErr markers (grammar-checker artifact for the double space).
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Build-Package([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
}

$q13 = $d.Paragraphs(13)
$q13Body = '<w:p>' +
             '<w:r><w:t xml:space="preserve">Q13) Find missing number </w:t></w:r>' +
             '<w:proofErr w:type="gramStart"/>' +
             '<w:r><w:t>in  a</w:t></w:r>' +
             '<w:proofErr w:type="gramEnd"/>' +
             '<w:r><w:t xml:space="preserve"> range from 1 to N</w:t></w:r>' +
           '</w:p>'
$q13.Range.InsertXML((Build-Package $q13Body))

# ---------------------------------------------------------------------------
# Q30) split "array(" out into its own run wrapped in gramStart/gramEnd
#      proofErr markers.
# ---------------------------------------------------------------------------
$q30 = $d.Paragraphs(30)
$q30Body = '<w:p>' +
             '<w:r><w:t xml:space="preserve">Q30) Find majority element in an </w:t></w:r>' +
             '<w:proofErr w:type="gramStart"/>' +
             '<w:r><w:t>array(</w:t></w:r>' +
             '<w:proofErr w:type="gramEnd"/>' +
             '<w:r><w:t>if any)</w:t></w:r>' +
           '</w:p>'
$q30.Range.InsertXML((Build-Package $q30Body))
